$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# coinranking.com crypto price/volume refresh (scheduled GitHub Actions run).
# Price (D) and Volume(1h) (E) are stored as plain text in this sheet (values like
# "24.569.09" use "." as a thousands separator and "  +3.50%  " keeps its padding),
# so each cell is pinned to Text just long enough to write the literal string -- this
# stops Excel from "helpfully" reinterpreting it as a number/date and mangling it (e.g.
# dropping a trailing zero) -- then the style is put back to Normal immediately after.
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "24.569.09"
Set-TextValue "E2" "  +3.50%  "

Set-TextValue "D3" "1.693.77"
Set-TextValue "E3" "  +1.77%  "

Set-TextValue "D4" "1.001"
Set-TextValue "E4" "  +0.12%  "

Set-TextValue "D5" "315.61"
Set-TextValue "E5" "  +1.99%  "

Set-TextValue "D7" "0.3934"
Set-TextValue "E7" "  +1.44%  "

Set-TextValue "D8" "0.4009"
Set-TextValue "E8" "  +1.83%  "

Set-TextValue "D9" "1.526"
Set-TextValue "E9" "  +5.94%  "

Set-TextValue "D10" "1.001"
Set-TextValue "E10" "  +0.13%  "

Set-TextValue "D11" "52.86"
Set-TextValue "E11" "  +6.35%  "

Set-TextValue "D12" "0.08735"
Set-TextValue "E12" "  +0.98%  "

Set-TextValue "E13" "  +7.31%  "

Set-TextValue "E14" "  +2.40%  "

Set-TextValue "D15" "0.00001316"
Set-TextValue "E15" "  +0.39%  "

Set-TextValue "D16" "7.564"
Set-TextValue "E16" "  +4.34%  "

Set-TextValue "D17" "1.693.99"
Set-TextValue "E17" "  +1.69%  "

Set-TextValue "D18" "99.69"
Set-TextValue "E18" "  +0.22%  "

Set-TextValue "D19" "0.07052"
Set-TextValue "E19" "  +4.02%  "

Set-TextValue "E20" "  +3.37%  "

Set-TextValue "D21" "6.859"
Set-TextValue "E21" "  +3.52%  "

Set-TextValue "E22" "  +0.02%  "

Set-TextValue "D23" "14.04"
Set-TextValue "E23" "  +1.70%  "

Set-TextValue "D24" "24.567.48"
Set-TextValue "E24" "  +3.52%  "

Set-TextValue "D25" "3.006"
Set-TextValue "E25" "  +6.86%  "

Set-TextValue "D26" "2.319"
Set-TextValue "E26" "  -0.12%  "

Set-TextValue "E27" "  +2.99%  "

Set-TextValue "D28" "160.07"
Set-TextValue "E28" "  +0.59%  "

Set-TextValue "D29" "5.217"
Set-TextValue "E29" "  +0.97%  "

Set-TextValue "D30" "134.45"
Set-TextValue "E30" "  +3.89%  "

Set-TextValue "D31" "7.454"
Set-TextValue "E31" "  +10.51%  "

Set-TextValue "D32" "1.881.31"
Set-TextValue "E32" "  +1.81%  "

Set-TextValue "D33" "1.097"
Set-TextValue "E33" "  -1.95%  "

Set-TextValue "D34" "0.08505"

Set-TextValue "D35" "7.214"
Set-TextValue "E35" "  +8.73%  "

Set-TextValue "D36" "11.47"
Set-TextValue "E36" "  +9.72%  "

Set-TextValue "D37" "1.961"
Set-TextValue "E37" "  +0.11%  "

Set-TextValue "E38" "  +2.58%  "

Set-TextValue "E39" "  -0.16%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D40" "0.02741"
Set-TextValue "E40" "  +9.25%  "

$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D41" "0.09043"
Set-TextValue "E41" "  +3.08%  "

Set-TextValue "D42" "1.472"
Set-TextValue "E42" "  +1.50%  "

Set-TextValue "E43" "  +2.26%  "

Set-TextValue "D44" "0.7187"
Set-TextValue "E44" "  +2.71%  "

$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D45" "2.537"
Set-TextValue "E45" "  +5.63%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D46" "15.33"
Set-TextValue "E46" "  +3.69%  "

Set-TextValue "D47" "4.211"
Set-TextValue "E47" "  +2.75%  "

Set-TextValue "E48" "  +0.14%  "

Set-TextValue "D49" "140.88"
Set-TextValue "E49" "  +1.65%  "

Set-TextValue "D50" "1.321"
Set-TextValue "E50" "  +7.31%  "

Set-TextValue "D51" "0.08010"
Set-TextValue "E51" "  +3.23%  "
